# doc/开发进度.xlsx — add properties for EEllipse in QStudioSCADA and QSCADARunTime.
#
# 椭圆 (Ellipse) row moves from "进行中" (in progress) to "已完成" (completed),
# and 直线 (Line) row moves from "未开始" (not started) to "进行中" (in progress).
# Refresh the status column's color coding to match: green = 已完成,
# yellow = 进行中, red = 未开始.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- status text updates -------------------------------------------------
$ws.Range("B5").Value = "已完成"   # 椭圆 (Ellipse): 进行中 -> 已完成
$ws.Range("B6").Value = "进行中"   # 直线 (Line): 未开始 -> 进行中

# --- status color coding ---------------------------------------------------
# Green fill (00B050) for completed rows
$ws.Range("B2:B5").Interior.Color = 0x50B000
# Yellow fill (FFFF00) for in-progress rows
$ws.Range("B6").Interior.Color = 0x00FFFF
# Red fill (FF0000) for not-started rows
$ws.Range("B7:B12").Interior.Color = 0x0000FF

# --- page setup --------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- selection -----------------------------------------------------------
$ws.Range("B19").Select()
